$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 562
$ws.Range("F5").Value = 2477
$ws.Range("F7").Value = 147
$ws.Range("F8").Value = 12
$ws.Range("F9").Value = 234
$ws.Range("F10").Value = 4970
$ws.Range("F11").Value = 6236
$ws.Range("F12").Value = 890
$ws.Range("F13").Value = 84
$ws.Range("F14").Value = 1396
$ws.Range("F15").Value = 1342
$ws.Range("F16").Value = 552
$ws.Range("F17").Value = 6794
$ws.Range("F18").Value = 378
$ws.Range("F19").Value = 28
$ws.Range("F20").Value = 62
$ws.Range("F21").Value = 4570
$ws.Range("F22").Value = 375
$ws.Range("F23").Value = 29
$ws.Range("F24").Value = 747
$ws.Range("F25").Value = 2225
$ws.Range("F26").Value = 1230
$ws.Range("F27").Value = 415
$ws.Range("F28").Value = 1130
$ws.Range("F29").Value = 177
$ws.Range("F30").Value = 78
$ws.Range("F31").Value = 64
$ws.Range("F32").Value = 130
$ws.Range("F33").Value = 358
$ws.Range("F34").Value = 1249
$ws.Range("F35").Value = 1956
$ws.Range("F36").Value = 197
$ws.Range("F37").Value = 481
$ws.Range("F39").Value = 1318
$ws.Range("F42").Value = 3
$ws.Range("F43").Value = 83
$ws.Range("F44").Value = 1055
$ws.Range("F45").Value = 1341
$ws.Range("F47").Value = 88
$ws.Range("F48").Value = 217
$ws.Range("F49").Value = 58

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 423
$ws.Range("F5").Value = 448
$ws.Range("F6").Value = 237
$ws.Range("F10").Value = 3
$ws.Range("F11").Value = 2
$ws.Range("F13").Value = 242
$ws.Range("F15").Value = 22
$ws.Range("F20").Value = 124
$ws.Range("F21").Value = 15
$ws.Range("F26").Value = 301
$ws.Range("F27").Value = 252
$ws.Range("F37").Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1625
$ws.Range("F7").Value = 524
$ws.Range("F8").Value = 3186
$ws.Range("F9").Value = 1204
$ws.Range("F10").Value = 1196
$ws.Range("F12").Value = 2045
$ws.Range("F13").Value = 505
$ws.Range("F14").Value = 417

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1625
$ws.Range("F3").Value = 562
$ws.Range("F4").Value = 524
$ws.Range("F5").Value = 2477
$ws.Range("F6").Value = 1204
$ws.Range("F7").Value = 234
$ws.Range("F8").Value = 2045
$ws.Range("F9").Value = 4970
$ws.Range("F10").Value = 505
$ws.Range("F11").Value = 448
$ws.Range("F12").Value = 237
$ws.Range("F13").Value = 890
$ws.Range("F14").Value = 84
$ws.Range("F16").Value = 1396
$ws.Range("F17").Value = 1342
$ws.Range("F18").Value = 552
$ws.Range("F19").Value = 6794
$ws.Range("F20").Value = 378
$ws.Range("F21").Value = 417
$ws.Range("F23").Value = 3
$ws.Range("F24").Value = 4570
$ws.Range("F25").Value = 375
$ws.Range("F26").Value = 747
$ws.Range("F27").Value = 2225
$ws.Range("F28").Value = 1230
$ws.Range("F29").Value = 415
$ws.Range("F30").Value = 1130
$ws.Range("F31").Value = 177
$ws.Range("F32").Value = 64
$ws.Range("F33").Value = 242
$ws.Range("F35").Value = 130
$ws.Range("F36").Value = 358
$ws.Range("F37").Value = 1956
$ws.Range("F38").Value = 197
$ws.Range("F39").Value = 481
$ws.Range("F41").Value = 15
$ws.Range("F42").Value = 1318
$ws.Range("F45").Value = 301
$ws.Range("F47").Value = 1055
$ws.Range("F48").Value = 1341
$ws.Range("F49").Value = 217
